$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header labels ---
$ws.Range("A1").Value = "button_apiKeyManagement_trNthChild"
$ws.Range("B1").Value = "button_apiKeyManagement_trNthChild_1"
$ws.Range("C1").Value = "button_apiKeyManagement_trNthChild_2"
$ws.Range("D1").Value = "input_KeyName"

# Give the new header cells the same ("Pandas") formatting that A1 already
# carries, by copying A1's format onto them (keeps them on the same cell
# style record instead of allocating new duplicate style entries).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 2 values (stored as text, not numbers) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"
$ws.Range("C2").Style = "Normal"

# D2 stays empty but present in the sheet data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = ""
$ws.Range("D2").Style = "Normal"

# --- Column widths (stored width = ColumnWidth + 5/6, so subtract that
#     back off to land exactly on the target stored widths of 36/38/38/15) ---
$ws.Columns.Item(1).ColumnWidth = 36 - 5/6
$ws.Columns.Item(2).ColumnWidth = 38 - 5/6
$ws.Columns.Item(3).ColumnWidth = 38 - 5/6
$ws.Columns.Item(4).ColumnWidth = 15 - 5/6
